$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6899999999999999
$ws.Range("C2").Value = 0.01
$ws.Range("D2").Value = 2.9
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.42
$ws.Range("G2").Value = 2.38
$ws.Range("H2").Value = 1.59
$ws.Range("I2").Value = 0.99
$ws.Range("J2").Value = 0.92
$ws.Range("K2").Value = 0.73

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.42
$ws.Range("G3").Value = 2.95
$ws.Range("H3").Value = 2.98
$ws.Range("I3").Value = 1.96
$ws.Range("J3").Value = 0.9
$ws.Range("K3").Value = 0.74

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.42
$ws.Range("G4").Value = 1.21
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 0.48
$ws.Range("J4").Value = 0.9
$ws.Range("K4").Value = 0.73
